# Allow pdf output to use extracted report type and notes
#
# This updates the label/placeholder text on both quarterly-report sheets so
# that downstream pdf-export code can pull the report type and notes that
# were actually extracted, instead of the old generic "screened / not
# screened" wording. It also mirrors the author's view-state changes
# (which sheet/cell is selected, and the widened label column).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# ---------------------------------------------------------------------
# Sheet1 ("previous quarter") text updates
# ---------------------------------------------------------------------
$ws1.Range("B1").Value = 'Name of report type (e.g. "Mammogram screening")'
$ws1.Range("B2").Value = "Notes describing total patients"
$ws1.Range("B3").Value = "Notes describing criteria for patients"

$ws1.Range("A10").Value = "Patients meeting criteria"
$ws1.Range("A11").Value = "Patients not meeting criteria"
$ws1.Range("A17").Value = "Patients meeting criteria"
$ws1.Range("A18").Value = "Patients not meeting criteria"
$ws1.Range("A26").Value = "Patients meeting criteria"
$ws1.Range("A27").Value = "Patients not meeting criteria"

# Widen the label column to fit the new, longer text
$ws1.Columns.Item(1).ColumnWidth = 39.57142857142857

# ---------------------------------------------------------------------
# Sheet2 ("current quarter") text updates
# ---------------------------------------------------------------------
$ws2.Range("B1").Value = 'Name of report type (e.g. "Mammogram screening")'
$ws2.Range("B2").Value = "Notes describing total patients"
$ws2.Range("B3").Value = "Notes describing criteria for patients"

$ws2.Range("A10").Value = "Patients meeting criteria"
$ws2.Range("A11").Value = "Patients not meeting criteria"
$ws2.Range("A17").Value = "Patients meeting criteria"
$ws2.Range("A18").Value = "Patients not meeting criteria"

$ws2.Columns.Item(1).ColumnWidth = 39.285714285714285

# ---------------------------------------------------------------------
# View-state: select G5 on Sheet1, then make Sheet2 the active sheet
# with F10 selected (matches the saved workbook view in the commit)
# ---------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("G5").Select()

$ws2.Activate()
$ws2.Range("F10").Select()
